$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 288.1579
$ws.Range("I33").Value = 286.66666
$ws.Range("J33").Value = 298
$ws.Range("K33").Value = 286.66666
$ws.Range("L33").Value = 298
$ws.Range("M33").Value = -57.66665999999998
$ws.Range("N33").Value = -756
$ws.Range("H88").Value = 2065.4333
$ws.Range("I88").Value = 1231.6666
$ws.Range("J88").Value = 2273.875
$ws.Range("K88").Value = 1231.6666
$ws.Range("L88").Value = 2273.875
$ws.Range("M88").Value = -825.6666
$ws.Range("N88").Value = -3085.875
$ws.Range("H91").Value = 2065.4333
$ws.Range("I91").Value = 1231.6666
$ws.Range("J91").Value = 2273.875
$ws.Range("K91").Value = 1231.6666
$ws.Range("L91").Value = 2273.875
$ws.Range("M91").Value = 172.3334
$ws.Range("N91").Value = -5081.875
$ws.Range("H92").Value = 493.1
$ws.Range("I92").Value = 505.2857
$ws.Range("J92").Value = 464.66666
$ws.Range("K92").Value = 505.2857
$ws.Range("L92").Value = 464.66666
$ws.Range("M92").Value = 742.7143
$ws.Range("N92").Value = -2960.66666
$ws.Range("H100").Value = 1435
$ws.Range("I100").Value = 1180
$ws.Range("J100").Value = 1690
$ws.Range("K100").Value = 1180
$ws.Range("L100").Value = 1690
$ws.Range("M100").Value = -639
$ws.Range("N100").Value = -2772
$ws.Range("H103").Value = 1121.3667
$ws.Range("I103").Value = 918
$ws.Range("J103").Value = 1256.9445
$ws.Range("K103").Value = 2754
$ws.Range("L103").Value = 3770.8335
$ws.Range("M103").Value = -2168
$ws.Range("N103").Value = -4942.833500000001
$ws.Range("H106").Value = 3500
$ws.Range("I106").Value = 2500
$ws.Range("K106").Value = 2500
$ws.Range("M106").Value = -1869
$ws.Range("H129").Value = 1054.826
$ws.Range("J129").Value = 1191.2368
$ws.Range("L129").Value = 3573.7104
$ws.Range("N129").Value = -13573.7104
$ws.Range("H132").Value = 198446.83
$ws.Range("I132").Value = 229862.7
$ws.Range("J132").Value = 975.5714
$ws.Range("K132").Value = 689588.1000000001
$ws.Range("L132").Value = 2926.7142
$ws.Range("M132").Value = -687058.1000000001
$ws.Range("N132").Value = -7986.7142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 2998.2
$ws.Range("I15").Value = 2998.2
$ws.Range("K15").Value = 2998.2
$ws.Range("M15").Value = -2648.2
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H61").Value = 2338.348
$ws.Range("I61").Value = 1992.6364
$ws.Range("J61").Value = 3215.923
$ws.Range("K61").Value = 1992.6364
$ws.Range("L61").Value = 3215.923
$ws.Range("M61").Value = -1780.6364
$ws.Range("N61").Value = -3639.923
$ws.Range("H74").Value = 2384135.5
$ws.Range("I74").Value = 3030791
$ws.Range("K74").Value = 3030791
$ws.Range("M74").Value = -3029917
$ws.Range("H77").Value = 2384135.5
$ws.Range("I77").Value = 3030791
$ws.Range("K77").Value = 15153955
$ws.Range("M77").Value = -15149587
$ws.Range("H97").Value = 362.32352
$ws.Range("I97").Value = 314.84616
$ws.Range("J97").Value = 516.625
$ws.Range("K97").Value = 314.84616
$ws.Range("L97").Value = 516.625
$ws.Range("M97").Value = 181.15384
$ws.Range("N97").Value = -1508.625
$ws.Range("H132").Value = 16290501
$ws.Range("I132").Value = 29863750
$ws.Range("J132").Value = 2601.2
$ws.Range("K132").Value = 89591250
$ws.Range("L132").Value = 7803.599999999999
$ws.Range("M132").Value = -89588720
$ws.Range("N132").Value = -12863.6
$ws.Range("H136").Value = 2338.348
$ws.Range("I136").Value = 1992.6364
$ws.Range("J136").Value = 3215.923
$ws.Range("K136").Value = 5977.9092
$ws.Range("L136").Value = 9647.769
$ws.Range("M136").Value = -3427.9092
$ws.Range("N136").Value = -14747.769

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2192.2307
$ws.Range("I20").Value = 2799.8
$ws.Range("J20").Value = 1812.5
$ws.Range("K20").Value = 2799.8
$ws.Range("L20").Value = 1812.5
$ws.Range("M20").Value = -2552.8
$ws.Range("N20").Value = -2306.5
$ws.Range("H22").Value = 211.22223
$ws.Range("I22").Value = 211.22223
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 211.22223
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -38.22223
$ws.Range("N22").ClearContents()
$ws.Range("H134").Value = 56750.35
$ws.Range("I134").Value = 85546.69500000001
$ws.Range("J134").Value = 3271.4285
$ws.Range("K134").Value = 256640.085
$ws.Range("L134").Value = 9814.2855
$ws.Range("M134").Value = -254105.085
$ws.Range("N134").Value = -14884.2855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1927.8923
$ws.Range("I58").Value = 794.85297
$ws.Range("J58").Value = 3170.5806
$ws.Range("K58").Value = 794.85297
$ws.Range("L58").Value = 3170.5806
$ws.Range("M58").Value = -591.85297
$ws.Range("N58").Value = -3576.5806
$ws.Range("H132").Value = 4320.75
$ws.Range("I132").Value = 4655.8184
$ws.Range("K132").Value = 13967.4552
$ws.Range("M132").Value = -11437.4552
$ws.Range("H136").Value = 1927.8923
$ws.Range("I136").Value = 794.85297
$ws.Range("J136").Value = 3170.5806
$ws.Range("K136").Value = 2384.55891
$ws.Range("L136").Value = 9511.7418
$ws.Range("M136").Value = 165.4410899999998
$ws.Range("N136").Value = -14611.7418

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1024.8125
$ws.Range("I68").Value = 678.1539
$ws.Range("J68").Value = 1434.5
$ws.Range("K68").Value = 2034.4617
$ws.Range("L68").Value = 4303.5
$ws.Range("M68").Value = -1223.4617
$ws.Range("N68").Value = -5925.5
$ws.Range("H71").Value = 1024.8125
$ws.Range("I71").Value = 678.1539
$ws.Range("J71").Value = 1434.5
$ws.Range("K71").Value = 6103.3851
$ws.Range("L71").Value = 12910.5
$ws.Range("M71").Value = -2047.3851
$ws.Range("N71").Value = -21022.5
$ws.Range("H82").Value = 4333.3335
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2594
$ws.Range("H85").Value = 4333.3335
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1596
$ws.Range("H129").Value = 797.15
$ws.Range("I129").Value = 272.5
$ws.Range("J129").Value = 1146.9166
$ws.Range("K129").Value = 817.5
$ws.Range("L129").Value = 3440.7498
$ws.Range("M129").Value = 4182.5
$ws.Range("N129").Value = -13440.7498
$ws.Range("H131").Value = 1589260.8
$ws.Range("J131").Value = 1925181.4
$ws.Range("L131").Value = 5775544.199999999
$ws.Range("N131").Value = -5785624.199999999
$ws.Range("H132").Value = 100001380
$ws.Range("I132").Value = 166668380
$ws.Range("J132").Value = 862.25
$ws.Range("K132").Value = 1500015420
$ws.Range("L132").Value = 7760.25
$ws.Range("M132").Value = -1500012890
$ws.Range("N132").Value = -12820.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H122").Value = 4440.952
$ws.Range("I122").Value = 4440.952
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13322.856
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10872.856
$ws.Range("N122").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7417.3184
$ws.Range("I122").Value = 11487.272
$ws.Range("J122").Value = 3347.3635
$ws.Range("K122").Value = 34461.81600000001
$ws.Range("L122").Value = 10042.0905
$ws.Range("M122").Value = -32011.81600000001
$ws.Range("N122").Value = -14942.0905
$ws.Range("H136").Value = 2008.25
$ws.Range("I136").Value = 1568.8518
$ws.Range("J136").Value = 4381
$ws.Range("K136").Value = 4706.555399999999
$ws.Range("L136").Value = 13143
$ws.Range("M136").Value = -2156.555399999999
$ws.Range("N136").Value = -18243

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2227.182
$ws.Range("I81").Value = 1820
$ws.Range("J81").Value = 2566.5
$ws.Range("K81").Value = 3640
$ws.Range("L81").Value = 5133
$ws.Range("M81").Value = -2579
$ws.Range("N81").Value = -7255
$ws.Range("H84").Value = 2227.182
$ws.Range("I84").Value = 1820
$ws.Range("J84").Value = 2566.5
$ws.Range("K84").Value = 18200
$ws.Range("L84").Value = 25665
$ws.Range("M84").Value = -12896
$ws.Range("N84").Value = -36273
$ws.Range("H96").Value = 8343426
$ws.Range("J96").Value = 13178.889
$ws.Range("L96").Value = 13178.889
$ws.Range("N96").Value = -15924.889
$ws.Range("H122").Value = 1581.9546
$ws.Range("I122").Value = 1378.4286
$ws.Range("J122").Value = 1938.125
$ws.Range("K122").Value = 4135.2858
$ws.Range("L122").Value = 5814.375
$ws.Range("M122").Value = -1685.2858
$ws.Range("N122").Value = -10714.375
$ws.Range("H132").Value = 9125217
$ws.Range("I132").Value = 11210482
$ws.Range("J132").Value = 2180.375
$ws.Range("K132").Value = 33631446
$ws.Range("L132").Value = 6541.125
$ws.Range("M132").Value = -33628916
$ws.Range("N132").Value = -11601.125
$ws.Range("H136").Value = 6004134.5
$ws.Range("I136").Value = 18025.5
$ws.Range("J136").Value = 11990244
$ws.Range("K136").Value = 54076.5
$ws.Range("L136").Value = 35970732
$ws.Range("M136").Value = -51526.5
$ws.Range("N136").Value = -35975832

Write-Host "Applied all market data updates."